$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.872.95'
$ws.Cells.Item(2, 5).Value = '  -0.35%  '
$ws.Cells.Item(3, 4).Value = '2.364.57'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.82'
$ws.Cells.Item(5, 5).Value = '  -0.05%  '
$ws.Cells.Item(6, 2).Value = 'XRP'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.670'
$ws.Cells.Item(6, 5).Value = '  -1.83%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '74.18'
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.599'
$ws.Cells.Item(9, 5).Value = '  +1.00%  '
$ws.Cells.Item(10, 5).Value = '  +2.00%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '60.12'
$ws.Cells.Item(11, 5).Value = '  +4.94%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.18'
$ws.Cells.Item(12, 5).Value = '  +15.48%  '
$ws.Cells.Item(13, 5).Value = '  +0.58%  '
$ws.Cells.Item(14, 5).Value = '  +0.30%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.44'
$ws.Cells.Item(15, 5).Value = '  -0.62%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.927'
$ws.Cells.Item(16, 5).Value = '  +2.79%  '
$ws.Cells.Item(17, 4).Value = '2.371.02'
$ws.Cells.Item(17, 5).Value = '  +0.47%  '
$ws.Cells.Item(18, 4).Value = '43.832.09'
$ws.Cells.Item(18, 5).Value = '  -0.26%  '
$ws.Cells.Item(19, 5).Value = '  +2.33%  '
$ws.Cells.Item(20, 5).Value = '  -5.81%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '77.29'
$ws.Cells.Item(21, 5).Value = '  -0.10%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '253.96'
$ws.Cells.Item(22, 5).Value = '  -1.89%  '
$ws.Cells.Item(23, 5).Value = '  +3.46%  '
$ws.Cells.Item(24, 5).Value = '  +0.03%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.86'
$ws.Cells.Item(25, 5).Value = '  -5.38%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.50'
$ws.Cells.Item(26, 5).Value = '  +0.40%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.64'
$ws.Cells.Item(27, 5).Value = '  -1.53%  '
$ws.Cells.Item(28, 5).Value = '  +0.73%  '
$ws.Cells.Item(29, 5).Value = '  -1.95%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '175.13'
$ws.Cells.Item(30, 5).Value = '  -0.21%  '
$ws.Cells.Item(31, 5).Value = '  +0.63%  '
$ws.Cells.Item(32, 5).Value = '  -2.09%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0756'
$ws.Cells.Item(33, 5).Value = '  -0.31%  '
$ws.Cells.Item(34, 5).Value = '  -2.49%  '
$ws.Cells.Item(35, 5).Value = '  -2.08%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.81'
$ws.Cells.Item(36, 5).Value = '  +1.10%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.62'
$ws.Cells.Item(37, 5).Value = '  +4.19%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0279'
$ws.Cells.Item(39, 5).Value = '  -0.03%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.65'
$ws.Cells.Item(40, 5).Value = '  +19.13%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.73'
$ws.Cells.Item(41, 5).Value = '  +9.74%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.98'
$ws.Cells.Item(42, 5).Value = '  +9.31%  '
$ws.Cells.Item(43, 5).Value = '  -4.23%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.07'
$ws.Cells.Item(44, 5).Value = '  +0.78%  '
$ws.Cells.Item(45, 5).Value = '  -1.05%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.52'
$ws.Cells.Item(46, 5).Value = '  +1.14%  '
$ws.Cells.Item(47, 2).Value = 'BinanceUSD'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$ws.Cells.Item(47, 5).Value = '  -0.05%  '
$ws.Cells.Item(48, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.24'
$ws.Cells.Item(48, 5).Value = '  -0.04%  '
$ws.Cells.Item(49, 5).Value = '  -0.60%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '98.46'
$ws.Cells.Item(50, 5).Value = '  -2.74%  '
$ws.Cells.Item(51, 5).Value = '  +2.22%  '
